$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Log Time column (D2:D106): new timestamp value, time-only
# number format, explicit black font color.
$timeRange = $ws.Range("D2:D106")
$timeRange.Value = 0.460590277777778
$timeRange.NumberFormat = "h:mm:ss"
$timeRange.Font.Color = 0

# Rows no longer carry an explicit height override.
$ws.Rows("1:106").AutoFit()

# Selection moves onto the edited column.
$timeRange.Select()
